$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9-14 down to 10-15
$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = "NARE "
$ws.Cells.Item(9, 2).Value = "N"
$ws.Cells.Item(9, 3).Value = "JAZMIN,ABARCO,GIRASOL"

$ws.Range("C10").Select()
